# Update automatico via Actualizar 04-16-2021 13-47-43
# Refreshes the "last updated" timestamp column (D) for each of the three
# time-block groups on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D15").Value = 44302.57454448682
$ws.Range("D16:D29").Value = 44302.55307650463
$ws.Range("D30:D43").Value = 44302.53166145834
